$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New SIP entries for 2025-01-01 (Excel serial date 45658), appended below
# the existing table (which currently ends at row 7).
$rows = @(
    @(7,  "Parag Parikh Flexi cap",        "0P0000YWL1.BO", 45658, "SIP", 1500),
    @(8,  "Nippon India Small cap",        "0P0000XVFY.BO", 45658, "SIP", 1000),
    @(9,  "DSP Nifty 50 Equal Weightage",  "0P0001BOXZ.BO", 45658, "SIP", 500),
    @(10, "DSP Nifty Next 50",             "0P0001FTFQ.BO", 45658, "SIP", 500),
    @(11, "Motilal Oswal Nasdaq 100 FoF",  "0P0001F0CK.BO", 45658, "SIP", 500)
)

$destRow = 8
foreach ($row in $rows) {
    # Clone the formatting of the row directly above so the new row's
    # cell styles (number formats, alignment, etc.) match the existing
    # table without introducing any new style/numFmt entries.
    $ws.Range("A" + ($destRow - 1) + ":F" + ($destRow - 1)).Copy()
    $ws.Range("A" + $destRow + ":F" + $destRow).PasteSpecial(-4122)

    $ws.Cells.Item($destRow, 1).Value = $row[0]
    $ws.Cells.Item($destRow, 2).Value = $row[1]
    $ws.Cells.Item($destRow, 3).Value = $row[2]
    $ws.Cells.Item($destRow, 4).Value = $row[3]
    $ws.Cells.Item($destRow, 5).Value = $row[4]
    $ws.Cells.Item($destRow, 6).Value = $row[5]

    $destRow++
}

$ws.Application.CutCopyMode = $false
$ws.Range("C12").Select()
